$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, -0.2359795975604194, 1.900377010157683, 11.8897934829752, 3.448157984051079, 3.521027694101534, 22),
    @(3, -0.3407653498485861, 1.038808833104293, 4.33724652853252, 2.08260570644866, 2.105274751693671, 21),
    @(4, -0.131555750701869, 0.8460533655601525, 2.692507442697656, 1.640886176033443, 1.678094322408583, 20),
    @(5, 0.0945621327354168, 0.8356365543419514, 2.101234218072367, 1.449563457759738, 1.486112606047758, 19),
    @(6, -0.05968872511646385, 0.6371319026568236, 0.778561022117543, 0.8823610497509186, 0.9058622479295383, 18),
    @(7, 0.09265814781738935, 0.5989980499374025, 0.6944883350010943, 0.8333596672512381, 0.8536812799923659, 17),
    @(8, 0.189786657949447, 0.5225847098032543, 0.4703199534050894, 0.6857987703438155, 0.6806279280150794, 16),
    @(9, 0.2020970918481325, 0.4802119434926625, 0.4448954675707786, 0.6670048482363368, 0.6579613262386209, 15),
    @(10, 0.216926569107512, 0.3902862058645244, 0.2053025722297817, 0.4531032688359043, 0.4128173466307566, 14),
    @(11, 0.2444166901726743, 0.3999205401577337, 0.2150394398791692, 0.4637234519400212, 0.4101726241042355, 13)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
